# Updated symbol list on Wed Dec 28 22:40:11 UTC 2022 with GitHub Actions
#
# Applies the cryptos.xlsx price/row updates described by the commit diff:
#   - Column D ("Price") cells are refreshed with new quotes. These must stay
#     TEXT (the sheet stores prices as inline strings, trailing zeros and
#     all -- e.g. "243.10", "23.70"), so a plain .Value assignment of a
#     numeric-looking string would be auto-coerced into a real number by
#     Excel and lose the exact text formatting. Prefixing with a leading
#     apostrophe forces text entry (like typing '243.10 into the cell),
#     and re-applying the "Normal" style afterwards clears the resulting
#     quote-prefix style flag so no stray formatting is introduced.
#   - Columns B/C/E hold plain text (coin name, link, volume label) and are
#     assigned directly.
#   - Rows 18-24 are a full rank rotation: each row's old Coin/Link/Volume
#     trio shifts up one rank (row 19's identity moves into row 18, etc.),
#     row 24 wraps around to what used to be row 18's coin ("One"), and
#     every row also gets a freshly scraped Price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice {
    param($addr, $value)
    # Leading apostrophe = "treat as text" (mirrors typing '<value> in Excel),
    # which keeps numeric-looking strings (trailing zeros etc.) intact.
    $ws.Range($addr).Value = "'" + $value
    # Clear the quote-prefix formatting flag that the apostrophe entry adds,
    # so no unintended style change is left behind.
    $ws.Range($addr).Style = "Normal"
}

# -- Row 2..17: price-only refreshes -----------------------------------
Set-TextPrice "D2"  "243.10"
Set-TextPrice "D3"  "23.70"
Set-TextPrice "D4"  "5.233"
Set-TextPrice "D5"  "0.05756"
Set-TextPrice "D6"  "6.406"
Set-TextPrice "D7"  "3.227"
Set-TextPrice "D8"  "0.8051"
Set-TextPrice "D9"  "0.8869"
Set-TextPrice "D10" "0.1373"
Set-TextPrice "D11" "0.07074"
Set-TextPrice "D13" "0.03039"
Set-TextPrice "D14" "0.09302"
Set-TextPrice "D15" "3.808"
Set-TextPrice "D16" "0.001548"
Set-TextPrice "D17" "0.04707"

# -- Row 18: was One -> now TigerCash -----------------------------------
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextPrice "D18" "0.006196"
$ws.Range("E18").Value = "17TigerCashTCH"

# -- Row 19: was TigerCash -> now BitKan --------------------------------
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextPrice "D19" "0.001250"
$ws.Range("E19").Value = "18BitKanKAN"

# -- Row 20: was BitKan -> now HotbitToken ------------------------------
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextPrice "D20" "0.004060"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# -- Row 21: was HotbitToken -> now NitroEx -----------------------------
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextPrice "D21" "0.00008697"
$ws.Range("E21").Value = "20NitroExNTX"

# -- Row 22: was NitroEx -> now LEO --------------------------------------
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextPrice "D22" "3.545"
$ws.Range("E22").Value = "21LEOLEO"

# -- Row 23: was LEO -> now BTSEToken ------------------------------------
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextPrice "D23" "2.140"
$ws.Range("E23").Value = "22BTSETokenBTSE"

# -- Row 24: was BTSEToken -> now One (wraps back around) ----------------
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextPrice "D24" "0.0006028"
$ws.Range("E24").Value = "23OneONE"

# -- Remaining scattered price / label refreshes -------------------------
Set-TextPrice "D26" "0.1318"
Set-TextPrice "D28" "0.0002328"
Set-TextPrice "D40" "0.03721"

Set-TextPrice "D41" "0.006248"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

Set-TextPrice "D42" "0.1043"
Set-TextPrice "D43" "0.002470"
Set-TextPrice "D44" "0.007148"
Set-TextPrice "D45" "0.00005336"

Set-TextPrice "D47" "0.5348"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

Set-TextPrice "D48" "0.002417"
$ws.Range("E48").Value = "47BOLOBOLO"

Set-TextPrice "D49" "0.00002099"
Set-TextPrice "D50" "0.0001999"
